# Quarterly indexing esoteric bug-fix operation
#
# Column A (rows 2-63) holds the "as-of" date that each row's naive AR
# forecast was produced from. These were previously stamped as the 1st of
# the month; the corrected indexing scheme re-stamps them on the 15th of
# the *following* month instead. Shift every date in A2:A63 accordingly,
# leaving styles/number-formats and every other column untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$epoch = Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row   # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $serial = $cell.Value2
    if ($null -ne $serial) {
        $dt = $epoch.AddDays($serial)

        $newMonth = $dt.Month + 1
        $newYear = $dt.Year
        if ($newMonth -gt 12) {
            $newMonth = 1
            $newYear = $newYear + 1
        }

        $newDate = Get-Date -Year $newYear -Month $newMonth -Day 15 -Hour 0 -Minute 0 -Second 0
        $cell.Value = $newDate
    }
}
